$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("uMudd Mark V")
$ws2 = $wb.Worksheets.Item("Tester Board")

# --- Update the 0.1uF capacitor row: C17 added to the designator list, qty 5 -> 8 ---
$ws1.Range("A20").Value = "C4, C6, C8, C11, C12, C13, C15, C17"
$ws1.Range("G20").Value = 8

# --- New crystal (Y1) row replaces the old blank "total" row 27 ---
$ws1.Range("A27").Value = "Y1"
$ws1.Range("B27").Value = "40MHz oscillator"
$ws1.Range("C27").Value = "TXC Corporation"
$ws1.Range("D27").Value = "9B-40.000MAAJ-B"
$ws1.Range("E27").Value = "Digikey"
$ws1.Range("F27").Value = "887-2030-ND"
$ws1.Range("G27").Value = 1
$ws1.Range("H27").Value = 0.39
$ws1.Range("I27").Formula = "=G27*H27"
$ws1.Range("J27").Value = "through-hole"
$ws1.Range("K27").Value = "in stock"

# --- Grand-total formula moves down from I27 to I30 ---
$ws1.Range("I30").Formula = "=SUM(I3:I27)"

# --- Column A widened to fit the longer designator list ---
# (target stored width 32.453125; the host quantizes ColumnWidth to 1/6-character
# steps, so 31.65 is the closest settable value that lands on that grid at 32.5)
$ws1.Columns.Item(1).ColumnWidth = 31.65

# --- Selection / active-tab bookkeeping: uMudd Mark V becomes the active tab ---
$ws2.Activate()
$ws2.Range("A18").Select()

$ws1.Activate()
$ws1.Range("I31").Select()
